# Auto-generated edit script: swaps the two match rows that were
# recorded in reversed order for 12 fixture dates, and refreshes four
# closing-odds cells on the final (still-unplayed) fixture row 339.
# Mirrors the Mexico Liga MX base update commit (22-05-2024 20:16).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap data row 34 <-> data row 35 (columns B, E:AB; A/C/D stay put)
$ws.Cells.Item(34,2).Value2 = 7094230
$ws.Cells.Item(35,2).Value2 = 7053868
$ws.Cells.Item(34,5).Value2 = "Club America"
$ws.Cells.Item(35,5).Value2 = "Necaxa"
$ws.Cells.Item(34,6).Value2 = "Atlas"
$ws.Cells.Item(35,6).Value2 = "Tigres UANL"
$ws.Cells.Item(34,7).Value2 = 1
$ws.Cells.Item(35,7).Value2 = 0
$ws.Cells.Item(34,8).Value2 = 1
$ws.Cells.Item(35,8).Value2 = 3
$ws.Cells.Item(34,9).Value2 = "D"
$ws.Cells.Item(35,9).Value2 = "A"
$ws.Cells.Item(34,10).Value2 = 2.3
$ws.Cells.Item(35,10).Value2 = 3.75
$ws.Cells.Item(34,11).Value2 = 3.4
$ws.Cells.Item(35,11).Value2 = 3.3
$ws.Cells.Item(34,12).Value2 = 3.1
$ws.Cells.Item(35,12).Value2 = 1.909
$ws.Cells.Item(34,13).Value2 = 2.15
$ws.Cells.Item(35,13).Value2 = 4.5
$ws.Cells.Item(34,14).Value2 = 3.4
$ws.Cells.Item(35,14).Value2 = 3.6
$ws.Cells.Item(34,15).Value2 = 3.4
$ws.Cells.Item(35,15).Value2 = 1.833
$ws.Cells.Item(34,16).Value2 = -0.25
$ws.Cells.Item(35,16).Value2 = 0.5
$ws.Cells.Item(34,17).Value2 = 1.85
$ws.Cells.Item(35,17).Value2 = 2.025
$ws.Cells.Item(34,18).Value2 = 2
$ws.Cells.Item(35,18).Value2 = 1.825
$ws.Cells.Item(34,19).Value2 = 2.75
$ws.Cells.Item(35,19).Value2 = 2.5
$ws.Cells.Item(34,20).Value2 = 2
$ws.Cells.Item(35,20).Value2 = 2
$ws.Cells.Item(34,21).Value2 = 1.85
$ws.Cells.Item(35,21).Value2 = 1.85
$ws.Cells.Item(34,22).Value2 = -1
$ws.Cells.Item(35,22).Value2 = -1
$ws.Cells.Item(34,23).Value2 = 2.4
$ws.Cells.Item(35,23).Value2 = -1
$ws.Cells.Item(34,24).Value2 = -1
$ws.Cells.Item(35,24).Value2 = 0.833
$ws.Cells.Item(34,25).Value2 = -0.5
$ws.Cells.Item(35,25).Value2 = -1
$ws.Cells.Item(34,26).Value2 = 0.5
$ws.Cells.Item(35,26).Value2 = 0.825
$ws.Cells.Item(34,27).Value2 = -1
$ws.Cells.Item(35,27).Value2 = 1
$ws.Cells.Item(34,28).Value2 = 0.8500000000000001
$ws.Cells.Item(35,28).Value2 = -1

# Swap data row 72 <-> data row 73 (columns B, E:AB; A/C/D stay put)
$ws.Cells.Item(72,2).Value2 = 6754048
$ws.Cells.Item(73,2).Value2 = 6754049
$ws.Cells.Item(72,5).Value2 = "Atletico San Luis"
$ws.Cells.Item(73,5).Value2 = "Juarez FC"
$ws.Cells.Item(72,6).Value2 = "Mazatlan FC"
$ws.Cells.Item(73,6).Value2 = "Atlas"
$ws.Cells.Item(72,7).Value2 = 3
$ws.Cells.Item(73,7).Value2 = 1
$ws.Cells.Item(72,8).Value2 = 2
$ws.Cells.Item(73,8).Value2 = 2
$ws.Cells.Item(72,9).Value2 = "H"
$ws.Cells.Item(73,9).Value2 = "A"
$ws.Cells.Item(72,10).Value2 = 1.615
$ws.Cells.Item(73,10).Value2 = 2.75
$ws.Cells.Item(72,11).Value2 = 4
$ws.Cells.Item(73,11).Value2 = 3.25
$ws.Cells.Item(72,12).Value2 = 4.5
$ws.Cells.Item(73,12).Value2 = 2.375
$ws.Cells.Item(72,13).Value2 = 1.6
$ws.Cells.Item(73,13).Value2 = 2.6
$ws.Cells.Item(72,14).Value2 = 4.5
$ws.Cells.Item(73,14).Value2 = 3.2
$ws.Cells.Item(72,15).Value2 = 5
$ws.Cells.Item(73,15).Value2 = 2.8
$ws.Cells.Item(72,16).Value2 = -1
$ws.Cells.Item(73,16).Value2 = 0
$ws.Cells.Item(72,17).Value2 = 1.95
$ws.Cells.Item(73,17).Value2 = 1.85
$ws.Cells.Item(72,18).Value2 = 1.9
$ws.Cells.Item(73,18).Value2 = 2
$ws.Cells.Item(72,19).Value2 = 3
$ws.Cells.Item(73,19).Value2 = 2.25
$ws.Cells.Item(72,20).Value2 = 1.925
$ws.Cells.Item(73,20).Value2 = 2.1
$ws.Cells.Item(72,21).Value2 = 1.925
$ws.Cells.Item(73,21).Value2 = 1.775
$ws.Cells.Item(72,22).Value2 = 0.6000000000000001
$ws.Cells.Item(73,22).Value2 = -1
$ws.Cells.Item(72,23).Value2 = -1
$ws.Cells.Item(73,23).Value2 = -1
$ws.Cells.Item(72,24).Value2 = -1
$ws.Cells.Item(73,24).Value2 = 1.8
$ws.Cells.Item(72,25).Value2 = 0
$ws.Cells.Item(73,25).Value2 = -1
$ws.Cells.Item(72,26).Value2 = 0
$ws.Cells.Item(73,26).Value2 = 1
$ws.Cells.Item(72,27).Value2 = 0.925
$ws.Cells.Item(73,27).Value2 = 1.1
$ws.Cells.Item(72,28).Value2 = -1
$ws.Cells.Item(73,28).Value2 = -1

# Swap data row 92 <-> data row 93 (columns B, E:AB; A/C/D stay put)
$ws.Cells.Item(92,2).Value2 = 6754065
$ws.Cells.Item(93,2).Value2 = 6754066
$ws.Cells.Item(92,5).Value2 = "Necaxa"
$ws.Cells.Item(93,5).Value2 = "Unam Pumas"
$ws.Cells.Item(92,6).Value2 = "Cruz Azul"
$ws.Cells.Item(93,6).Value2 = "Queretaro"
$ws.Cells.Item(92,7).Value2 = 1
$ws.Cells.Item(93,7).Value2 = 4
$ws.Cells.Item(92,8).Value2 = 3
$ws.Cells.Item(93,8).Value2 = 0
$ws.Cells.Item(92,9).Value2 = "A"
$ws.Cells.Item(93,9).Value2 = "H"
$ws.Cells.Item(92,10).Value2 = 2.375
$ws.Cells.Item(93,10).Value2 = 1.727
$ws.Cells.Item(92,11).Value2 = 3.3
$ws.Cells.Item(93,11).Value2 = 3.5
$ws.Cells.Item(92,12).Value2 = 2.8
$ws.Cells.Item(93,12).Value2 = 4.5
$ws.Cells.Item(92,13).Value2 = 3.5
$ws.Cells.Item(93,13).Value2 = 1.8
$ws.Cells.Item(92,14).Value2 = 3.6
$ws.Cells.Item(93,14).Value2 = 3.6
$ws.Cells.Item(92,15).Value2 = 2.1
$ws.Cells.Item(93,15).Value2 = 4.5
$ws.Cells.Item(92,16).Value2 = 0.25
$ws.Cells.Item(93,16).Value2 = -0.75
$ws.Cells.Item(92,17).Value2 = 2
$ws.Cells.Item(93,17).Value2 = 2.025
$ws.Cells.Item(92,18).Value2 = 1.85
$ws.Cells.Item(93,18).Value2 = 1.825
$ws.Cells.Item(92,19).Value2 = 2.5
$ws.Cells.Item(93,19).Value2 = 2.5
$ws.Cells.Item(92,20).Value2 = 1.9
$ws.Cells.Item(93,20).Value2 = 1.825
$ws.Cells.Item(92,21).Value2 = 1.95
$ws.Cells.Item(93,21).Value2 = 2.025
$ws.Cells.Item(92,22).Value2 = -1
$ws.Cells.Item(93,22).Value2 = 0.8
$ws.Cells.Item(92,23).Value2 = -1
$ws.Cells.Item(93,23).Value2 = -1
$ws.Cells.Item(92,24).Value2 = 1.1
$ws.Cells.Item(93,24).Value2 = -1
$ws.Cells.Item(92,25).Value2 = -1
$ws.Cells.Item(93,25).Value2 = 1.025
$ws.Cells.Item(92,26).Value2 = 0.8500000000000001
$ws.Cells.Item(93,26).Value2 = -1
$ws.Cells.Item(92,27).Value2 = 0.8999999999999999
$ws.Cells.Item(93,27).Value2 = 0.825
$ws.Cells.Item(92,28).Value2 = -1
$ws.Cells.Item(93,28).Value2 = -1

# Swap data row 98 <-> data row 99 (columns B, E:AB; A/C/D stay put)
$ws.Cells.Item(98,2).Value2 = 6754641
$ws.Cells.Item(99,2).Value2 = 6754074
$ws.Cells.Item(98,5).Value2 = "Pachuca"
$ws.Cells.Item(99,5).Value2 = "Chivas Guadalajara"
$ws.Cells.Item(98,6).Value2 = "Tigres UANL"
$ws.Cells.Item(99,6).Value2 = "Atlas"
$ws.Cells.Item(98,7).Value2 = 1
$ws.Cells.Item(99,7).Value2 = 4
$ws.Cells.Item(98,8).Value2 = 1
$ws.Cells.Item(99,8).Value2 = 1
$ws.Cells.Item(98,9).Value2 = "D"
$ws.Cells.Item(99,9).Value2 = "H"
$ws.Cells.Item(98,10).Value2 = 2.875
$ws.Cells.Item(99,10).Value2 = 2.3
$ws.Cells.Item(98,11).Value2 = 3.5
$ws.Cells.Item(99,11).Value2 = 3.3
$ws.Cells.Item(98,12).Value2 = 2.2
$ws.Cells.Item(99,12).Value2 = 2.8
$ws.Cells.Item(98,13).Value2 = 2.9
$ws.Cells.Item(99,13).Value2 = 2.4
$ws.Cells.Item(98,14).Value2 = 3.5
$ws.Cells.Item(99,14).Value2 = 3.2
$ws.Cells.Item(98,15).Value2 = 2.4
$ws.Cells.Item(99,15).Value2 = 3.1
$ws.Cells.Item(98,16).Value2 = 0.25
$ws.Cells.Item(99,16).Value2 = -0.25
$ws.Cells.Item(98,17).Value2 = 1.775
$ws.Cells.Item(99,17).Value2 = 2.15
$ws.Cells.Item(98,18).Value2 = 2.1
$ws.Cells.Item(99,18).Value2 = 1.725
$ws.Cells.Item(98,19).Value2 = 2.5
$ws.Cells.Item(99,19).Value2 = 2
$ws.Cells.Item(98,20).Value2 = 1.825
$ws.Cells.Item(99,20).Value2 = 1.925
$ws.Cells.Item(98,21).Value2 = 2.025
$ws.Cells.Item(99,21).Value2 = 1.925
$ws.Cells.Item(98,22).Value2 = -1
$ws.Cells.Item(99,22).Value2 = 1.4
$ws.Cells.Item(98,23).Value2 = 2.5
$ws.Cells.Item(99,23).Value2 = -1
$ws.Cells.Item(98,24).Value2 = -1
$ws.Cells.Item(99,24).Value2 = -1
$ws.Cells.Item(98,25).Value2 = 0.3875
$ws.Cells.Item(99,25).Value2 = 1.15
$ws.Cells.Item(98,26).Value2 = -0.5
$ws.Cells.Item(99,26).Value2 = -1
$ws.Cells.Item(98,27).Value2 = -1
$ws.Cells.Item(99,27).Value2 = 0.925
$ws.Cells.Item(98,28).Value2 = 1.025
$ws.Cells.Item(99,28).Value2 = -1

# Swap data row 132 <-> data row 133 (columns B, E:AB; A/C/D stay put)
$ws.Cells.Item(132,2).Value2 = 6754129
$ws.Cells.Item(133,2).Value2 = 6754103
$ws.Cells.Item(132,5).Value2 = "Atletico San Luis"
$ws.Cells.Item(133,5).Value2 = "Mazatlan FC"
$ws.Cells.Item(132,6).Value2 = "Club America"
$ws.Cells.Item(133,6).Value2 = "Santos Laguna"
$ws.Cells.Item(132,7).Value2 = 0
$ws.Cells.Item(133,7).Value2 = 3
$ws.Cells.Item(132,8).Value2 = 1
$ws.Cells.Item(133,8).Value2 = 1
$ws.Cells.Item(132,9).Value2 = "A"
$ws.Cells.Item(133,9).Value2 = "H"
$ws.Cells.Item(132,10).Value2 = 3.75
$ws.Cells.Item(133,10).Value2 = 2.5
$ws.Cells.Item(132,11).Value2 = 3.3
$ws.Cells.Item(133,11).Value2 = 3.3
$ws.Cells.Item(132,12).Value2 = 2
$ws.Cells.Item(133,12).Value2 = 2.75
$ws.Cells.Item(132,13).Value2 = 5
$ws.Cells.Item(133,13).Value2 = 2.8
$ws.Cells.Item(132,14).Value2 = 3.8
$ws.Cells.Item(133,14).Value2 = 3.4
$ws.Cells.Item(132,15).Value2 = 1.7
$ws.Cells.Item(133,15).Value2 = 2.5
$ws.Cells.Item(132,16).Value2 = 0.75
$ws.Cells.Item(133,16).Value2 = 0
$ws.Cells.Item(132,17).Value2 = 1.95
$ws.Cells.Item(133,17).Value2 = 2.05
$ws.Cells.Item(132,18).Value2 = 1.9
$ws.Cells.Item(133,18).Value2 = 1.8
$ws.Cells.Item(132,19).Value2 = 3
$ws.Cells.Item(133,19).Value2 = 3
$ws.Cells.Item(132,20).Value2 = 1.925
$ws.Cells.Item(133,20).Value2 = 1.975
$ws.Cells.Item(132,21).Value2 = 1.925
$ws.Cells.Item(133,21).Value2 = 1.875
$ws.Cells.Item(132,22).Value2 = -1
$ws.Cells.Item(133,22).Value2 = 1.8
$ws.Cells.Item(132,23).Value2 = -1
$ws.Cells.Item(133,23).Value2 = -1
$ws.Cells.Item(132,24).Value2 = 0.7
$ws.Cells.Item(133,24).Value2 = -1
$ws.Cells.Item(132,25).Value2 = -0.5
$ws.Cells.Item(133,25).Value2 = 1.05
$ws.Cells.Item(132,26).Value2 = 0.45
$ws.Cells.Item(133,26).Value2 = -1
$ws.Cells.Item(132,27).Value2 = -1
$ws.Cells.Item(133,27).Value2 = 0.9750000000000001
$ws.Cells.Item(132,28).Value2 = 0.925
$ws.Cells.Item(133,28).Value2 = -1

# Swap data row 175 <-> data row 176 (columns B, E:AB; A/C/D stay put)
$ws.Cells.Item(175,2).Value2 = 7612685
$ws.Cells.Item(176,2).Value2 = 7612675
$ws.Cells.Item(175,5).Value2 = "Tijuana"
$ws.Cells.Item(176,5).Value2 = "Monterrey"
$ws.Cells.Item(175,6).Value2 = "Club America"
$ws.Cells.Item(176,6).Value2 = "Puebla"
$ws.Cells.Item(175,7).Value2 = 0
$ws.Cells.Item(176,7).Value2 = 2
$ws.Cells.Item(175,8).Value2 = 2
$ws.Cells.Item(176,8).Value2 = 0
$ws.Cells.Item(175,9).Value2 = "A"
$ws.Cells.Item(176,9).Value2 = "H"
$ws.Cells.Item(175,10).Value2 = 4.5
$ws.Cells.Item(176,10).Value2 = 1.571
$ws.Cells.Item(175,11).Value2 = 3.5
$ws.Cells.Item(176,11).Value2 = 3.75
$ws.Cells.Item(175,12).Value2 = 1.8
$ws.Cells.Item(176,12).Value2 = 6
$ws.Cells.Item(175,13).Value2 = 2.45
$ws.Cells.Item(176,13).Value2 = 1.363
$ws.Cells.Item(175,14).Value2 = 3.3
$ws.Cells.Item(176,14).Value2 = 5
$ws.Cells.Item(175,15).Value2 = 2.875
$ws.Cells.Item(176,15).Value2 = 7.5
$ws.Cells.Item(175,16).Value2 = -0.25
$ws.Cells.Item(176,16).Value2 = -1.5
$ws.Cells.Item(175,17).Value2 = 2.1
$ws.Cells.Item(176,17).Value2 = 2.025
$ws.Cells.Item(175,18).Value2 = 1.775
$ws.Cells.Item(176,18).Value2 = 1.825
$ws.Cells.Item(175,19).Value2 = 2.25
$ws.Cells.Item(176,19).Value2 = 3
$ws.Cells.Item(175,20).Value2 = 1.8
$ws.Cells.Item(176,20).Value2 = 2
$ws.Cells.Item(175,21).Value2 = 2.05
$ws.Cells.Item(176,21).Value2 = 1.85
$ws.Cells.Item(175,22).Value2 = -1
$ws.Cells.Item(176,22).Value2 = 0.363
$ws.Cells.Item(175,23).Value2 = -1
$ws.Cells.Item(176,23).Value2 = -1
$ws.Cells.Item(175,24).Value2 = 1.875
$ws.Cells.Item(176,24).Value2 = -1
$ws.Cells.Item(175,25).Value2 = -1
$ws.Cells.Item(176,25).Value2 = 1.025
$ws.Cells.Item(175,26).Value2 = 0.7749999999999999
$ws.Cells.Item(176,26).Value2 = -1
$ws.Cells.Item(175,27).Value2 = -0.5
$ws.Cells.Item(176,27).Value2 = -1
$ws.Cells.Item(175,28).Value2 = 0.5249999999999999
$ws.Cells.Item(176,28).Value2 = 0.8500000000000001

# Swap data row 193 <-> data row 194 (columns B, E:AB; A/C/D stay put)
$ws.Cells.Item(193,2).Value2 = 7612811
$ws.Cells.Item(194,2).Value2 = 7612810
$ws.Cells.Item(193,5).Value2 = "Leon"
$ws.Cells.Item(194,5).Value2 = "Cruz Azul"
$ws.Cells.Item(193,6).Value2 = "Santos Laguna"
$ws.Cells.Item(194,6).Value2 = "Mazatlan FC"
$ws.Cells.Item(193,7).Value2 = 3
$ws.Cells.Item(194,7).Value2 = 2
$ws.Cells.Item(193,8).Value2 = 2
$ws.Cells.Item(194,8).Value2 = 1
$ws.Cells.Item(193,9).Value2 = "H"
$ws.Cells.Item(194,9).Value2 = "H"
$ws.Cells.Item(193,10).Value2 = 1.833
$ws.Cells.Item(194,10).Value2 = 1.727
$ws.Cells.Item(193,11).Value2 = 3.75
$ws.Cells.Item(194,11).Value2 = 3.9
$ws.Cells.Item(193,12).Value2 = 4
$ws.Cells.Item(194,12).Value2 = 4.333
$ws.Cells.Item(193,13).Value2 = 1.8
$ws.Cells.Item(194,13).Value2 = 1.5
$ws.Cells.Item(193,14).Value2 = 4
$ws.Cells.Item(194,14).Value2 = 4.2
$ws.Cells.Item(193,15).Value2 = 4.2
$ws.Cells.Item(194,15).Value2 = 5.5
$ws.Cells.Item(193,16).Value2 = -0.75
$ws.Cells.Item(194,16).Value2 = -1
$ws.Cells.Item(193,17).Value2 = 2
$ws.Cells.Item(194,17).Value2 = 1.85
$ws.Cells.Item(193,18).Value2 = 1.85
$ws.Cells.Item(194,18).Value2 = 2
$ws.Cells.Item(193,19).Value2 = 3
$ws.Cells.Item(194,19).Value2 = 2.75
$ws.Cells.Item(193,20).Value2 = 1.975
$ws.Cells.Item(194,20).Value2 = 1.8
$ws.Cells.Item(193,21).Value2 = 1.875
$ws.Cells.Item(194,21).Value2 = 2.05
$ws.Cells.Item(193,22).Value2 = 0.8
$ws.Cells.Item(194,22).Value2 = 0.5
$ws.Cells.Item(193,23).Value2 = -1
$ws.Cells.Item(194,23).Value2 = -1
$ws.Cells.Item(193,24).Value2 = -1
$ws.Cells.Item(194,24).Value2 = -1
$ws.Cells.Item(193,25).Value2 = 0.5
$ws.Cells.Item(194,25).Value2 = 0
$ws.Cells.Item(193,26).Value2 = -0.5
$ws.Cells.Item(194,26).Value2 = 0
$ws.Cells.Item(193,27).Value2 = 0.9750000000000001
$ws.Cells.Item(194,27).Value2 = 0.4
$ws.Cells.Item(193,28).Value2 = -1
$ws.Cells.Item(194,28).Value2 = -0.5

# Swap data row 237 <-> data row 238 (columns B, E:AB; A/C/D stay put)
$ws.Cells.Item(237,2).Value2 = 7612867
$ws.Cells.Item(238,2).Value2 = 7612866
$ws.Cells.Item(237,5).Value2 = "Club America"
$ws.Cells.Item(238,5).Value2 = "Leon"
$ws.Cells.Item(237,6).Value2 = "Mazatlan FC"
$ws.Cells.Item(238,6).Value2 = "Cruz Azul"
$ws.Cells.Item(237,7).Value2 = 2
$ws.Cells.Item(238,7).Value2 = 2
$ws.Cells.Item(237,8).Value2 = 2
$ws.Cells.Item(238,8).Value2 = 3
$ws.Cells.Item(237,9).Value2 = "D"
$ws.Cells.Item(238,9).Value2 = "A"
$ws.Cells.Item(237,10).Value2 = 1.363
$ws.Cells.Item(238,10).Value2 = 2.5
$ws.Cells.Item(237,11).Value2 = 5
$ws.Cells.Item(238,11).Value2 = 3.4
$ws.Cells.Item(237,12).Value2 = 7.5
$ws.Cells.Item(238,12).Value2 = 2.7
$ws.Cells.Item(237,13).Value2 = 1.222
$ws.Cells.Item(238,13).Value2 = 2.8
$ws.Cells.Item(237,14).Value2 = 6.5
$ws.Cells.Item(238,14).Value2 = 3.6
$ws.Cells.Item(237,15).Value2 = 12
$ws.Cells.Item(238,15).Value2 = 2.375
$ws.Cells.Item(237,16).Value2 = -1.75
$ws.Cells.Item(238,16).Value2 = 0.25
$ws.Cells.Item(237,17).Value2 = 1.825
$ws.Cells.Item(238,17).Value2 = 1.75
$ws.Cells.Item(237,18).Value2 = 2.025
$ws.Cells.Item(238,18).Value2 = 2.05
$ws.Cells.Item(237,19).Value2 = 3.25
$ws.Cells.Item(238,19).Value2 = 2.75
$ws.Cells.Item(237,20).Value2 = 1.975
$ws.Cells.Item(238,20).Value2 = 1.85
$ws.Cells.Item(237,21).Value2 = 1.875
$ws.Cells.Item(238,21).Value2 = 2
$ws.Cells.Item(237,22).Value2 = -1
$ws.Cells.Item(238,22).Value2 = -1
$ws.Cells.Item(237,23).Value2 = 5.5
$ws.Cells.Item(238,23).Value2 = -1
$ws.Cells.Item(237,24).Value2 = -1
$ws.Cells.Item(238,24).Value2 = 1.375
$ws.Cells.Item(237,25).Value2 = -1
$ws.Cells.Item(238,25).Value2 = -1
$ws.Cells.Item(237,26).Value2 = 1.025
$ws.Cells.Item(238,26).Value2 = 1.05
$ws.Cells.Item(237,27).Value2 = 0.9750000000000001
$ws.Cells.Item(238,27).Value2 = 0.8500000000000001
$ws.Cells.Item(237,28).Value2 = -1
$ws.Cells.Item(238,28).Value2 = -1

# Swap data row 272 <-> data row 273 (columns B, E:AB; A/C/D stay put)
$ws.Cells.Item(272,2).Value2 = 7612892
$ws.Cells.Item(273,2).Value2 = 7612894
$ws.Cells.Item(272,5).Value2 = "Leon"
$ws.Cells.Item(273,5).Value2 = "Tigres UANL"
$ws.Cells.Item(272,6).Value2 = "Puebla"
$ws.Cells.Item(273,6).Value2 = "Mazatlan FC"
$ws.Cells.Item(272,7).Value2 = 2
$ws.Cells.Item(273,7).Value2 = 5
$ws.Cells.Item(272,8).Value2 = 1
$ws.Cells.Item(273,8).Value2 = 1
$ws.Cells.Item(272,9).Value2 = "H"
$ws.Cells.Item(273,9).Value2 = "H"
$ws.Cells.Item(272,10).Value2 = 1.571
$ws.Cells.Item(273,10).Value2 = 1.4
$ws.Cells.Item(272,11).Value2 = 4
$ws.Cells.Item(273,11).Value2 = 4.5
$ws.Cells.Item(272,12).Value2 = 4.75
$ws.Cells.Item(273,12).Value2 = 6.5
$ws.Cells.Item(272,13).Value2 = 1.5
$ws.Cells.Item(273,13).Value2 = 1.615
$ws.Cells.Item(272,14).Value2 = 4.75
$ws.Cells.Item(273,14).Value2 = 4
$ws.Cells.Item(272,15).Value2 = 5.5
$ws.Cells.Item(273,15).Value2 = 5.5
$ws.Cells.Item(272,16).Value2 = -1
$ws.Cells.Item(273,16).Value2 = -1
$ws.Cells.Item(272,17).Value2 = 1.8
$ws.Cells.Item(273,17).Value2 = 2.05
$ws.Cells.Item(272,18).Value2 = 2.05
$ws.Cells.Item(273,18).Value2 = 1.8
$ws.Cells.Item(272,19).Value2 = 3.25
$ws.Cells.Item(273,19).Value2 = 2.75
$ws.Cells.Item(272,20).Value2 = 2.05
$ws.Cells.Item(273,20).Value2 = 1.925
$ws.Cells.Item(272,21).Value2 = 1.8
$ws.Cells.Item(273,21).Value2 = 1.925
$ws.Cells.Item(272,22).Value2 = 0.5
$ws.Cells.Item(273,22).Value2 = 0.615
$ws.Cells.Item(272,23).Value2 = -1
$ws.Cells.Item(273,23).Value2 = -1
$ws.Cells.Item(272,24).Value2 = -1
$ws.Cells.Item(273,24).Value2 = -1
$ws.Cells.Item(272,25).Value2 = 0
$ws.Cells.Item(273,25).Value2 = 1.05
$ws.Cells.Item(272,26).Value2 = 0
$ws.Cells.Item(273,26).Value2 = -1
$ws.Cells.Item(272,27).Value2 = -0.5
$ws.Cells.Item(273,27).Value2 = 0.925
$ws.Cells.Item(272,28).Value2 = 0.4
$ws.Cells.Item(273,28).Value2 = -1

# Swap data row 276 <-> data row 277 (columns B, E:AB; A/C/D stay put)
$ws.Cells.Item(276,2).Value2 = 7612897
$ws.Cells.Item(277,2).Value2 = 7612893
$ws.Cells.Item(276,5).Value2 = "Atletico San Luis"
$ws.Cells.Item(277,5).Value2 = "Toluca"
$ws.Cells.Item(276,6).Value2 = "Pachuca"
$ws.Cells.Item(277,6).Value2 = "Unam Pumas"
$ws.Cells.Item(276,7).Value2 = 2
$ws.Cells.Item(277,7).Value2 = 3
$ws.Cells.Item(276,8).Value2 = 1
$ws.Cells.Item(277,8).Value2 = 0
$ws.Cells.Item(276,9).Value2 = "H"
$ws.Cells.Item(277,9).Value2 = "H"
$ws.Cells.Item(276,10).Value2 = 2.55
$ws.Cells.Item(277,10).Value2 = 1.8
$ws.Cells.Item(276,11).Value2 = 3.75
$ws.Cells.Item(277,11).Value2 = 3.6
$ws.Cells.Item(276,12).Value2 = 2.3
$ws.Cells.Item(277,12).Value2 = 4
$ws.Cells.Item(276,13).Value2 = 3.5
$ws.Cells.Item(277,13).Value2 = 1.75
$ws.Cells.Item(276,14).Value2 = 3.6
$ws.Cells.Item(277,14).Value2 = 3.8
$ws.Cells.Item(276,15).Value2 = 2
$ws.Cells.Item(277,15).Value2 = 4.5
$ws.Cells.Item(276,16).Value2 = 0.5
$ws.Cells.Item(277,16).Value2 = -0.75
$ws.Cells.Item(276,17).Value2 = 1.825
$ws.Cells.Item(277,17).Value2 = 2
$ws.Cells.Item(276,18).Value2 = 2.025
$ws.Cells.Item(277,18).Value2 = 1.85
$ws.Cells.Item(276,19).Value2 = 2.75
$ws.Cells.Item(277,19).Value2 = 2.75
$ws.Cells.Item(276,20).Value2 = 1.8
$ws.Cells.Item(277,20).Value2 = 1.8
$ws.Cells.Item(276,21).Value2 = 2.05
$ws.Cells.Item(277,21).Value2 = 2.05
$ws.Cells.Item(276,22).Value2 = 2.5
$ws.Cells.Item(277,22).Value2 = 0.75
$ws.Cells.Item(276,23).Value2 = -1
$ws.Cells.Item(277,23).Value2 = -1
$ws.Cells.Item(276,24).Value2 = -1
$ws.Cells.Item(277,24).Value2 = -1
$ws.Cells.Item(276,25).Value2 = 0.825
$ws.Cells.Item(277,25).Value2 = 1
$ws.Cells.Item(276,26).Value2 = -1
$ws.Cells.Item(277,26).Value2 = -1
$ws.Cells.Item(276,27).Value2 = 0.4
$ws.Cells.Item(277,27).Value2 = 0.4
$ws.Cells.Item(276,28).Value2 = -0.5
$ws.Cells.Item(277,28).Value2 = -0.5

# Swap data row 298 <-> data row 299 (columns B, E:AB; A/C/D stay put)
$ws.Cells.Item(298,2).Value2 = 7612917
$ws.Cells.Item(299,2).Value2 = 7612918
$ws.Cells.Item(298,5).Value2 = "Necaxa"
$ws.Cells.Item(299,5).Value2 = "Puebla"
$ws.Cells.Item(298,6).Value2 = "Santos Laguna"
$ws.Cells.Item(299,6).Value2 = "Cruz Azul"
$ws.Cells.Item(298,7).Value2 = 2
$ws.Cells.Item(299,7).Value2 = 0
$ws.Cells.Item(298,8).Value2 = 0
$ws.Cells.Item(299,8).Value2 = 1
$ws.Cells.Item(298,9).Value2 = "H"
$ws.Cells.Item(299,9).Value2 = "A"
$ws.Cells.Item(298,10).Value2 = 2.4
$ws.Cells.Item(299,10).Value2 = 4.333
$ws.Cells.Item(298,11).Value2 = 3.25
$ws.Cells.Item(299,11).Value2 = 4.2
$ws.Cells.Item(298,12).Value2 = 2.875
$ws.Cells.Item(299,12).Value2 = 1.666
$ws.Cells.Item(298,13).Value2 = 2.05
$ws.Cells.Item(299,13).Value2 = 6.5
$ws.Cells.Item(298,14).Value2 = 3.3
$ws.Cells.Item(299,14).Value2 = 4.75
$ws.Cells.Item(298,15).Value2 = 3.8
$ws.Cells.Item(299,15).Value2 = 1.45
$ws.Cells.Item(298,16).Value2 = -0.5
$ws.Cells.Item(299,16).Value2 = 1.25
$ws.Cells.Item(298,17).Value2 = 2
$ws.Cells.Item(299,17).Value2 = 1.875
$ws.Cells.Item(298,18).Value2 = 1.85
$ws.Cells.Item(299,18).Value2 = 1.975
$ws.Cells.Item(298,19).Value2 = 2.5
$ws.Cells.Item(299,19).Value2 = 3.25
$ws.Cells.Item(298,20).Value2 = 1.95
$ws.Cells.Item(299,20).Value2 = 1.9
$ws.Cells.Item(298,21).Value2 = 1.9
$ws.Cells.Item(299,21).Value2 = 1.95
$ws.Cells.Item(298,22).Value2 = 1.05
$ws.Cells.Item(299,22).Value2 = -1
$ws.Cells.Item(298,23).Value2 = -1
$ws.Cells.Item(299,23).Value2 = -1
$ws.Cells.Item(298,24).Value2 = -1
$ws.Cells.Item(299,24).Value2 = 0.45
$ws.Cells.Item(298,25).Value2 = 1
$ws.Cells.Item(299,25).Value2 = 0.4375
$ws.Cells.Item(298,26).Value2 = -1
$ws.Cells.Item(299,26).Value2 = -0.5
$ws.Cells.Item(298,27).Value2 = -1
$ws.Cells.Item(299,27).Value2 = -1
$ws.Cells.Item(298,28).Value2 = 0.8999999999999999
$ws.Cells.Item(299,28).Value2 = 0.95

# Swap data row 322 <-> data row 323 (columns B, E:AB; A/C/D stay put)
$ws.Cells.Item(322,2).Value2 = 7612941
$ws.Cells.Item(323,2).Value2 = 8097226
$ws.Cells.Item(322,5).Value2 = "Necaxa"
$ws.Cells.Item(323,5).Value2 = "Santos Laguna"
$ws.Cells.Item(322,6).Value2 = "Monterrey"
$ws.Cells.Item(323,6).Value2 = "Atletico San Luis"
$ws.Cells.Item(322,7).Value2 = 2
$ws.Cells.Item(323,7).Value2 = 0
$ws.Cells.Item(322,8).Value2 = 5
$ws.Cells.Item(323,8).Value2 = 3
$ws.Cells.Item(322,9).Value2 = "A"
$ws.Cells.Item(323,9).Value2 = "A"
$ws.Cells.Item(322,10).Value2 = 3
$ws.Cells.Item(323,10).Value2 = 1.85
$ws.Cells.Item(322,11).Value2 = 3.5
$ws.Cells.Item(323,11).Value2 = 3.8
$ws.Cells.Item(322,12).Value2 = 2.25
$ws.Cells.Item(323,12).Value2 = 4
$ws.Cells.Item(322,13).Value2 = 3.25
$ws.Cells.Item(323,13).Value2 = 2.1
$ws.Cells.Item(322,14).Value2 = 3.4
$ws.Cells.Item(323,14).Value2 = 3.6
$ws.Cells.Item(322,15).Value2 = 2.2
$ws.Cells.Item(323,15).Value2 = 3.25
$ws.Cells.Item(322,16).Value2 = 0.25
$ws.Cells.Item(323,16).Value2 = -0.25
$ws.Cells.Item(322,17).Value2 = 1.975
$ws.Cells.Item(323,17).Value2 = 1.825
$ws.Cells.Item(322,18).Value2 = 1.875
$ws.Cells.Item(323,18).Value2 = 2.025
$ws.Cells.Item(322,19).Value2 = 2.5
$ws.Cells.Item(323,19).Value2 = 2.75
$ws.Cells.Item(322,20).Value2 = 1.875
$ws.Cells.Item(323,20).Value2 = 2
$ws.Cells.Item(322,21).Value2 = 1.975
$ws.Cells.Item(323,21).Value2 = 1.85
$ws.Cells.Item(322,22).Value2 = -1
$ws.Cells.Item(323,22).Value2 = -1
$ws.Cells.Item(322,23).Value2 = -1
$ws.Cells.Item(323,23).Value2 = -1
$ws.Cells.Item(322,24).Value2 = 1.2
$ws.Cells.Item(323,24).Value2 = 2.25
$ws.Cells.Item(322,25).Value2 = -1
$ws.Cells.Item(323,25).Value2 = -1
$ws.Cells.Item(322,26).Value2 = 0.875
$ws.Cells.Item(323,26).Value2 = 1.025
$ws.Cells.Item(322,27).Value2 = 0.875
$ws.Cells.Item(323,27).Value2 = 0.5
$ws.Cells.Item(322,28).Value2 = -1
$ws.Cells.Item(323,28).Value2 = -0.5

# Refresh closing Asian-handicap odds on row 339 (Cruz Azul vs Club America)
$ws.Cells.Item(339,17).Value2 = 1.95
$ws.Cells.Item(339,18).Value2 = 1.9
$ws.Cells.Item(339,20).Value2 = 2
$ws.Cells.Item(339,21).Value2 = 1.85

